$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 0.9872258258499474
$ws.Range("C3").Value = 0.9863239322320441
$ws.Range("D3").Value = 0.9862409830551231

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9849128769530381
$ws.Range("C4").Value = 0.9857990123122556
$ws.Range("D4").Value = 0.98592957943929

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9888030457699664
$ws.Range("C5").Value = 0.9889989083236905
$ws.Range("D5").Value = 0.9905120466406029
